$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apellidos / Nombres / No. Expediente (row 6)
$ws.Range("A6").Value = "GRAVE"
$ws.Range("C6").Value = "LOPEZ"
$ws.Range("E6").Value = "EULALIA"
$ws.Range("G6").Value = ""
$ws.Range("I6").Value = "2012-13786/201766035"

# Fecha de nacimiento / Edad / Lugar de nacimiento (row 12)
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1964-10-21"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "53"
$ws.Range("H12").Value = "BAJA VERAPAZ"

# Ocupacion / No. de Cedula (row 14)
$ws.Range("D14").Value = "AMA DE CASA"
$ws.Range("H14").Value = "NO PRESENTO"

# Contacto de emergencia: Nombre / Parentesco / Direccion / Telefono (row 20)
$ws.Range("A20").Value = "JOSE PILAR GRAVE"
$ws.Range("F20").Value = "HERMANO"
$ws.Range("H20").Value = "VILLA NUEVA"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "58102338"

# Fecha de ingreso / Hora (row 24)
$ws.Range("A24").Value = "20/11/2017"
$ws.Range("C24").Value = "15:13:51"
